$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 81 is the prior quarter's "AVERAGE of last 3" rollup row (date + AVERAGE
# formula). Append the next quarter (row 82) following the same pattern:
# date literal in column A, AVERAGE(prior 3 B cells) formula in column B.
$ws.Range("A82").Value = 44835
$ws.Range("B82").Formula = "=AVERAGE(B79:B81)"

# Carry the formatting (number format / font / borders / style) down from
# row 81 so the new row matches the existing rollup rows exactly.
$ws.Range("A81:B81").Copy()
$ws.Range("A82:B82").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection to the next empty row, mirroring where Excel leaves the
# cursor after entering data in A82/B82.
$ws.Range("A83").Select()

# Scroll the view down so the newly-entered row is visible (best effort —
# a no-op if the host doesn't track window scroll position).
$excel.ActiveWindow.ScrollRow = 70
